# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) - row -> new value for column F
$sheet1Updates = @{
    2  = 154
    3  = 1808
    4  = 40
    5  = 152
    6  = 674
    7  = 42
    9  = 560
    12 = 85
    13 = 170
    14 = 26
    18 = 5167
    20 = 845
    21 = 119
    22 = 2291
    23 = 72
    25 = 2138
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (sheetId 4) - row -> new value for column F
$sheet4Updates = @{
    2  = 154
    3  = 1808
    4  = 40
    5  = 152
    6  = 674
    7  = 42
    9  = 560
    12 = 85
    13 = 170
    14 = 26
    18 = 5167
    22 = 845
    23 = 119
    24 = 2291
    25 = 72
    28 = 2138
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
